$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing row 3 ACNO value (drop stray leading apostrophe character) ---
$ws.Range("A3").Value = "9943-000613-100"

# --- New "Currency" column values for the existing rows ---
$ws.Range("F2").Value = "US"
$ws.Range("F3").Value = "US"
$ws.Range("F4").Value = "US"

# --- Row 5: KTB / SINGAPORE / Saving / EUR account ---
$ws.Range("A5").Value = "9943-000613-002"
$ws.Range("B5").Value = "KTB"
$ws.Range("C5").Value = "SINGAPORE"
$ws.Range("D5").Value = "C.P.Trading"
$ws.Range("E5").Value = "Saving"
$ws.Range("F5").Value = "EUR"

# --- Row 6: CIMB / SINGAPORE / Saving / US account ---
$ws.Range("A6").Value = "2200027067340"
$ws.Range("B6").Value = "CIMB"
$ws.Range("C6").Value = "SINGAPORE"
$ws.Range("D6").Value = "C.P.Trading"
$ws.Range("E6").Value = "Saving"
$ws.Range("F6").Value = "US"

# --- Row 7: CIMB / SINGAPORE / Saving / EUR account ---
$ws.Range("A7").Value = "2000501927"
$ws.Range("B7").Value = "CIMB"
$ws.Range("C7").Value = "SINGAPORE"
$ws.Range("D7").Value = "C.P.Trading"
$ws.Range("E7").Value = "Saving"
$ws.Range("F7").Value = "EUR"

# --- Row 8: CIMB / SINGAPORE / Saving / CNH account ---
$ws.Range("A8").Value = "2000971203"
$ws.Range("B8").Value = "CIMB"
$ws.Range("C8").Value = "SINGAPORE"
$ws.Range("D8").Value = "C.P.Trading"
$ws.Range("E8").Value = "Saving"
$ws.Range("F8").Value = "CNH"

# --- Row 9: BBL -BKK (FCD) / BANGKOK / Saving / US account ---
$ws.Range("A9").Value = "840-101-0018-164202-001"
$ws.Range("B9").Value = "BBL -BKK (FCD)"
$ws.Range("C9").Value = "BANGKOK"
$ws.Range("D9").Value = "C.P.Trading"
$ws.Range("E9").Value = "Saving"
$ws.Range("F9").Value = "US"

# --- Grow the Table1 ListObject (and therefore sheet dimension + autofilter) to cover the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F9"))

# --- Column width tweaks (col A widened for the longer ACNO values, col B gets bestFit) ---
$ws.Columns.Item(1).ColumnWidth = 23.41666666666667
$ws.Columns.Item(2).ColumnWidth = 14.166666666666666

# --- Move the selection the way the author left it ---
$ws.Range("A9").Select()

# --- Restore the Excel window position recorded in the workbook view ---
$win = $wb.Windows.Item(1)
$win.Left = 6075
$win.Top = 6315
